$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New period columns BG:BI (31/12/2023, 31/03/2024, 30/06/2024), extending the
# existing quarterly balance-sheet layout that currently ends at column BF.

# Header row: copy the format of the last existing period header (BF1: bold,
# centered, bordered) onto the three new header cells, then set their labels.
$ws.Range("BF1").Copy()
$ws.Range("BG1:BI1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("BG1").Value = "31/12/2023"
$ws.Range("BH1").Value = "31/03/2024"
$ws.Range("BI1").Value = "30/06/2024"

$ws.Range("BG2").Value = 3076137.984
$ws.Range("BH2").Value = 2936739.072
$ws.Range("BI2").Value = 2962825.984
$ws.Range("BG3").Value = 138882
$ws.Range("BH3").Value = 107087
$ws.Range("BI3").Value = 161231.008
$ws.Range("BG4").Value = 17600
$ws.Range("BH4").Value = 14321
$ws.Range("BI4").Value = 78182
$ws.Range("BG5").Value = 5180
$ws.Range("BH5").Value = 4781
$ws.Range("BI5").Value = 4868
$ws.Range("BG6").Value = 32930
$ws.Range("BH6").Value = 26745
$ws.Range("BI6").Value = 30600
$ws.Range("BG7").Value = 0
$ws.Range("BH7").Value = 0
$ws.Range("BI7").Value = 0
$ws.Range("BG8").Value = 0
$ws.Range("BH8").Value = 0
$ws.Range("BI8").Value = 0
$ws.Range("BG9").Value = 20300
$ws.Range("BH9").Value = 17452
$ws.Range("BI9").Value = 12961
$ws.Range("BG10").Value = 1620
$ws.Range("BH10").Value = 1024
$ws.Range("BI10").Value = 1028
$ws.Range("BG11").Value = 61252
$ws.Range("BH11").Value = 42764
$ws.Range("BI11").Value = 33592
$ws.Range("BG12").Value = 315255.008
$ws.Range("BH12").Value = 226976
$ws.Range("BI12").Value = 217099.008
$ws.Range("BG13").Value = 0
$ws.Range("BH13").Value = 0
$ws.Range("BI13").Value = 0
$ws.Range("BG14").Value = 0
$ws.Range("BH14").Value = 0
$ws.Range("BI14").Value = 0
$ws.Range("BG15").Value = 312
$ws.Range("BH15").Value = 315
$ws.Range("BI15").Value = 21
$ws.Range("BG16").Value = 0
$ws.Range("BH16").Value = 0
$ws.Range("BI16").Value = 0
$ws.Range("BG17").Value = 0
$ws.Range("BH17").Value = 0
$ws.Range("BI17").Value = 0
$ws.Range("BG18").Value = 0
$ws.Range("BH18").Value = 0
$ws.Range("BI18").Value = 0
$ws.Range("BG19").Value = 4799
$ws.Range("BH19").Value = 5102
$ws.Range("BI19").Value = 3529
$ws.Range("BG20").Value = 0
$ws.Range("BH20").Value = 0
$ws.Range("BI20").Value = 0
$ws.Range("BG21").Value = 0
$ws.Range("BH21").Value = 0
$ws.Range("BI21").Value = 0
$ws.Range("BG22").Value = 0
$ws.Range("BH22").Value = 0
$ws.Range("BI22").Value = 0
$ws.Range("BG23").Value = 2622000.896
$ws.Range("BH23").Value = 2602675.968
$ws.Range("BI23").Value = 2584496.128
$ws.Range("BG24").Value = 0
$ws.Range("BH24").Value = 0
$ws.Range("BI24").Value = 0
$ws.Range("BG25").Value = 0
$ws.Range("BH25").Value = 0
$ws.Range("BI25").Value = 0
$ws.Range("BG26").Value = 3076137.984
$ws.Range("BH26").Value = 2936739.072
$ws.Range("BI26").Value = 2962825.984
$ws.Range("BG27").Value = 437913.984
$ws.Range("BH27").Value = 473134.016
$ws.Range("BI27").Value = 506687.008
$ws.Range("BG28").Value = 7220
$ws.Range("BH28").Value = 8766
$ws.Range("BI28").Value = 6101
$ws.Range("BG29").Value = 63750
$ws.Range("BH29").Value = 46843
$ws.Range("BI29").Value = 45795
$ws.Range("BG30").Value = 15033
$ws.Range("BH30").Value = 9511
$ws.Range("BI30").Value = 6807
$ws.Range("BG31").Value = 172754
$ws.Range("BH31").Value = 225503.008
$ws.Range("BI31").Value = 253178
$ws.Range("BG32").Value = 0
$ws.Range("BH32").Value = 0
$ws.Range("BI32").Value = 0
$ws.Range("BG33").Value = 0
$ws.Range("BH33").Value = 0
$ws.Range("BI33").Value = 0
$ws.Range("BG34").Value = 179156.992
$ws.Range("BH34").Value = 182511.008
$ws.Range("BI34").Value = 194806
$ws.Range("BG35").Value = 0
$ws.Range("BH35").Value = 0
$ws.Range("BI35").Value = 0
$ws.Range("BG36").Value = 0
$ws.Range("BH36").Value = 0
$ws.Range("BI36").Value = 0
$ws.Range("BG37").Value = 1705203.968
$ws.Range("BH37").Value = 1592500.992
$ws.Range("BI37").Value = 1571192.064
$ws.Range("BG38").Value = 901872
$ws.Range("BH38").Value = 876820.992
$ws.Range("BI38").Value = 876849.9840000001
$ws.Range("BG39").Value = 0
$ws.Range("BH39").Value = 0
$ws.Range("BI39").Value = 0
$ws.Range("BG40").Value = 417232
$ws.Range("BH40").Value = 407500
$ws.Range("BI40").Value = 428267.008
$ws.Range("BG41").Value = 0
$ws.Range("BH41").Value = 0
$ws.Range("BI41").Value = 0
$ws.Range("BG42").Value = 0
$ws.Range("BH42").Value = 0
$ws.Range("BI42").Value = 0
$ws.Range("BG43").Value = 386100
$ws.Range("BH43").Value = 308180
$ws.Range("BI43").Value = 266075.008
$ws.Range("BG44").Value = 0
$ws.Range("BH44").Value = 0
$ws.Range("BI44").Value = 0
$ws.Range("BG45").Value = 0
$ws.Range("BH45").Value = 0
$ws.Range("BI45").Value = 0
$ws.Range("BG46").Value = 0
$ws.Range("BH46").Value = 0
$ws.Range("BI46").Value = 0
$ws.Range("BG47").Value = 933020.032
$ws.Range("BH47").Value = 871104
$ws.Range("BI47").Value = 884947.008
$ws.Range("BG48").Value = 4128636.928
$ws.Range("BH48").Value = 4128636.928
$ws.Range("BI48").Value = 4128636.928
$ws.Range("BG49").Value = 1
$ws.Range("BH49").Value = 1
$ws.Range("BI49").Value = 1
$ws.Range("BG50").Value = 0
$ws.Range("BH50").Value = 0
$ws.Range("BI50").Value = 0
$ws.Range("BG51").Value = 0
$ws.Range("BH51").Value = 0
$ws.Range("BI51").Value = 0
$ws.Range("BG52").Value = -3195618.048
$ws.Range("BH52").Value = -3257533.952
$ws.Range("BI52").Value = -3243691.008
$ws.Range("BG53").Value = 0
$ws.Range("BH53").Value = 0
$ws.Range("BI53").Value = 0
$ws.Range("BG54").Value = 0
$ws.Range("BH54").Value = 0
$ws.Range("BI54").Value = 0
$ws.Range("BG55").Value = 0
$ws.Range("BH55").Value = 0
$ws.Range("BI55").Value = 0
$ws.Range("BG56").Value = 0
$ws.Range("BH56").Value = 0
$ws.Range("BI56").Value = 0
$ws.Range("BG57:BI57").Value = "'"
$ws.Range("BG57:BI57").ClearFormats()
$ws.Range("BG58:BI58").Value = "'"
$ws.Range("BG58:BI58").ClearFormats()
$ws.Range("BG59").Value = 50584
$ws.Range("BH59").Value = 42587
$ws.Range("BI59").Value = 56022
$ws.Range("BG60").Value = 104706
$ws.Range("BH60").Value = -57538
$ws.Range("BI60").Value = -59112
$ws.Range("BG61").Value = 155289.984
$ws.Range("BH61").Value = -14951
$ws.Range("BI61").Value = -3090
$ws.Range("BG62").Value = 0
$ws.Range("BH62").Value = 0
$ws.Range("BI62").Value = 0
$ws.Range("BG63").Value = -17934
$ws.Range("BH63").Value = -14421
$ws.Range("BI63").Value = -11615
$ws.Range("BG64").Value = 0
$ws.Range("BH64").Value = 0
$ws.Range("BI64").Value = 0
$ws.Range("BG65").Value = 0
$ws.Range("BH65").Value = 0
$ws.Range("BI65").Value = 0
$ws.Range("BG66").Value = 71126
$ws.Range("BH66").Value = -1058
$ws.Range("BI66").Value = 62834
$ws.Range("BG67").Value = 6878
$ws.Range("BH67").Value = 0
$ws.Range("BI67").Value = 0
$ws.Range("BG68").Value = -31035
$ws.Range("BH68").Value = -33814
$ws.Range("BI68").Value = -35947
$ws.Range("BG69").Value = 757
$ws.Range("BH69").Value = 512
$ws.Range("BI69").Value = 1378
$ws.Range("BG70").Value = -31792.008
$ws.Range("BH70").Value = -34326
$ws.Range("BI70").Value = -37325
$ws.Range("BG71:BI71").Value = "'"
$ws.Range("BG71:BI71").ClearFormats()
$ws.Range("BG72:BI72").Value = "'"
$ws.Range("BG72:BI72").ClearFormats()
$ws.Range("BG73:BI73").Value = "'"
$ws.Range("BG73:BI73").ClearFormats()
$ws.Range("BG74").Value = 184324.992
$ws.Range("BH74").Value = -64244
$ws.Range("BI74").Value = 12182
$ws.Range("BG75").Value = -2193
$ws.Range("BH75").Value = -1791
$ws.Range("BI75").Value = -16551
$ws.Range("BG76").Value = -56276
$ws.Range("BH76").Value = 4119
$ws.Range("BI76").Value = 18212
$ws.Range("BG77:BI77").Value = "'"
$ws.Range("BG77:BI77").ClearFormats()
$ws.Range("BG78:BI78").Value = "'"
$ws.Range("BG78:BI78").ClearFormats()
$ws.Range("BG79").Value = 0
$ws.Range("BH79").Value = 0
$ws.Range("BI79").Value = 0
$ws.Range("BG80").Value = 125856
$ws.Range("BH80").Value = -61916
$ws.Range("BI80").Value = 13843
